# Apply the changes described by the diff:
# - Brasil 2024 value changes from 4.36 to 4.43
# - A new "Brasil 2025" row (4.21) is inserted right after Brasil 2024
# - A new "Região Nordeste 2025" row (3.85) is inserted right after Região Nordeste 2024 (3.97, updated from 3.86)
# - A new "Sergipe 2025" row (3.72) is appended after Sergipe 2024, whose value changes from 4.32 to 4.73
# Net effect: each of the 3 regions (Brasil, Região Nordeste, Sergipe) now has rows for
# 2007..2025 (19 rows each), for a total of 57 data rows (dimension A1:C58).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to write a text value into column A or B while avoiding Excel's automatic
# conversion of date-shaped strings ("dd/mm/yyyy") into date serial numbers: a
# leading apostrophe forces the value to be stored as literal text.
function Set-TextCell($range, [string]$text) {
    $range.Value2 = "'" + $text
}

# 1) Update Brasil / 2024 value (row 19)
$ws.Range("C19").Value2 = 4.43

# 2) Insert a new row after row 19 for Brasil / 2025
$ws.Rows.Item(20).Insert()
Set-TextCell $ws.Range("A20") "Brasil"
Set-TextCell $ws.Range("B20") "01/01/2025"
$ws.Range("C20").Value2 = 4.21

# At this point the Região Nordeste block (originally rows 20-37) now sits at rows 21-38,
# with its last row (38) being year 2024. Update that value and insert the 2025 row after it.

# 3) Update Região Nordeste / 2024 value (now row 38)
$ws.Range("C38").Value2 = 3.97

# 4) Insert a new row after row 38 for Região Nordeste / 2025
$ws.Rows.Item(39).Insert()
Set-TextCell $ws.Range("A39") "Região Nordeste"
Set-TextCell $ws.Range("B39") "01/01/2025"
$ws.Range("C39").Value2 = 3.85

# The Sergipe block (originally rows 38-55) now sits at rows 40-57, with its last
# row (57) being year 2024.

# 5) Update Sergipe / 2024 value (now row 57)
$ws.Range("C57").Value2 = 4.73

# 6) Append a new row for Sergipe / 2025
Set-TextCell $ws.Range("A58") "Sergipe"
Set-TextCell $ws.Range("B58") "01/01/2025"
$ws.Range("C58").Value2 = 3.72
